# Apply quarterly financial data refresh for EBAY workbook
# 1) Insert two new columns (new quarter + prior quarter) before column D
# 2) Populate the new columns with data
# 3) Apply data corrections that came in with the refreshed quarterly figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two columns at D:E; existing D:K data shifts right to F:M
$ws.Range("D:E").Insert()

# Copy number formatting/styles from the (now shifted) data columns into the new D:E columns
$ws.Range("F7:M102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D (latest quarter) and column E (prior quarter) with data
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 2877000
$ws.Cells.Item(8, 5).Value = 2649000
$ws.Cells.Item(9, 4).Value = 618000
$ws.Cells.Item(9, 5).Value = 608000
$ws.Cells.Item(10, 4).Value = 2259000
$ws.Cells.Item(10, 5).Value = 2041000
$ws.Cells.Item(12, 4).Value = 292000
$ws.Cells.Item(12, 5).Value = 307000
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 13000
$ws.Cells.Item(15, 5).Value = 13000
$ws.Cells.Item(17, 4).Value = 2196000
$ws.Cells.Item(17, 5).Value = 2093000
$ws.Cells.Item(18, 4).Value = 681000
$ws.Cells.Item(18, 5).Value = 556000
$ws.Cells.Item(20, 4).Value = -165000
$ws.Cells.Item(20, 5).Value = 392000
$ws.Cells.Item(21, 4).Value = 688000
$ws.Cells.Item(21, 5).Value = 1123000
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 516000
$ws.Cells.Item(23, 5).Value = 948000
$ws.Cells.Item(24, 4).Value = 216000
$ws.Cells.Item(24, 5).Value = 228000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 300000
$ws.Cells.Item(26, 5).Value = 720000
$ws.Cells.Item(27, 4).Value = 300000
$ws.Cells.Item(27, 5).Value = 720000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 460000
$ws.Cells.Item(29, 5).Value = 1000
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 165000
$ws.Cells.Item(32, 5).Value = -392000
$ws.Cells.Item(33, 4).Value = 760000
$ws.Cells.Item(33, 5).Value = 721000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 760000
$ws.Cells.Item(35, 5).Value = 721000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 2202000
$ws.Cells.Item(41, 5).Value = 2086000
$ws.Cells.Item(42, 4).Value = 2696000
$ws.Cells.Item(42, 5).Value = 2737000
$ws.Cells.Item(43, 4).Value = 712000
$ws.Cells.Item(43, 5).Value = 761000
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(45, 4).Value = 1516000
$ws.Cells.Item(45, 5).Value = 1379000
$ws.Cells.Item(46, 4).Value = 7126000
$ws.Cells.Item(46, 5).Value = 6963000
$ws.Cells.Item(47, 4).Value = 3778000
$ws.Cells.Item(47, 5).Value = 4276000
$ws.Cells.Item(48, 4).Value = 1597000
$ws.Cells.Item(48, 5).Value = 1580000
$ws.Cells.Item(49, 4).Value = 5252000
$ws.Cells.Item(49, 5).Value = 5276000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 5066000
$ws.Cells.Item(52, 5).Value = 5557000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 22819000
$ws.Cells.Item(54, 5).Value = 23652000
$ws.Cells.Item(57, 4).Value = 286000
$ws.Cells.Item(57, 5).Value = 225000
$ws.Cells.Item(58, 4).Value = 1546000
$ws.Cells.Item(58, 5).Value = 1546000
$ws.Cells.Item(59, 4).Value = 2622000
$ws.Cells.Item(59, 5).Value = 2235000
$ws.Cells.Item(60, 4).Value = 4454000
$ws.Cells.Item(60, 5).Value = 4006000
$ws.Cells.Item(61, 4).Value = 7685000
$ws.Cells.Item(61, 5).Value = 7661000
$ws.Cells.Item(62, 4).Value = 4399000
$ws.Cells.Item(62, 5).Value = 5056000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 16538000
$ws.Cells.Item(66, 5).Value = 16723000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 16459000
$ws.Cells.Item(72, 5).Value = 15699000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 6281000
$ws.Cells.Item(76, 5).Value = 6929000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 760000
$ws.Cells.Item(81, 5).Value = 721000
$ws.Cells.Item(83, 4).Value = 172000
$ws.Cells.Item(83, 5).Value = 175000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 1233000
$ws.Cells.Item(89, 5).Value = 558000
$ws.Cells.Item(91, 4).Value = -130000
$ws.Cells.Item(91, 5).Value = -179000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = 421000
$ws.Cells.Item(94, 5).Value = 926000
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -1511000
$ws.Cells.Item(100, 5).Value = -1020000
$ws.Cells.Item(101, 4).Value = -25000
$ws.Cells.Item(101, 5).Value = 1000
$ws.Cells.Item(102, 4).Value = 118000
$ws.Cells.Item(102, 5).Value = 465000

# Apply corrections to previously-reported figures that shifted into columns F:M
$ws.Cells.Item(8, 8).Value = 2707000
$ws.Cells.Item(8, 9).Value = 2498000
$ws.Cells.Item(9, 9).Value = 557000
$ws.Cells.Item(10, 8).Value = 2117000
$ws.Cells.Item(10, 9).Value = 1941000
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(17, 8).Value = 2044000
$ws.Cells.Item(17, 9).Value = 1924000
$ws.Cells.Item(18, 9).Value = 574000
$ws.Cells.Item(20, 9).Value = 120000
$ws.Cells.Item(21, 9).Value = 867000
$ws.Cells.Item(23, 9).Value = 694000
$ws.Cells.Item(26, 9).Value = 520000
$ws.Cells.Item(27, 9).Value = 520000
$ws.Cells.Item(32, 9).Value = -120000
$ws.Cells.Item(33, 9).Value = 520000
$ws.Cells.Item(35, 9).Value = 520000
$ws.Cells.Item(81, 9).Value = 520000
$ws.Cells.Item(91, 9).Value = -157000
$ws.Cells.Item(91, 10).Value = -182000
